$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.019.73'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.682.14'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.07'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.58'
$ws.Range('E8').Value = '  +6.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.254'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0624'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = '1.921.13'
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('D13').Value = '1.653.95'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.29'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.21'
$ws.Range('E17').Value = '  +5.49%  '
$ws.Range('D18').Value = '27.042.86'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '236.27'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').Value = '0.0₃0738'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.27'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('E24').Value = '  -4.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.00'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.71'
$ws.Range('E26').Value = '  +5.01%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.24'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('E28').Value = '  -3.02%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0499'
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').Value = '1.523.22'
$ws.Range('E33').Value = '  +4.03%  '
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('E35').Value = '  +4.72%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.592'
$ws.Range('E37').Value = '  +3.63%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0176'
$ws.Range('E38').Value = '  +4.09%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.919'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.04'
$ws.Range('E40').Value = '  +7.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.76'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.07'
$ws.Range('E43').Value = '  +3.52%  '
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').Value = '1.825.00'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.781'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.28'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.104'
$ws.Range('E48').Value = '  +4.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.52'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.90'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0506'
